$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999999999954712
$ws.Range("E2").Value = 0.9999999999954712

# Row 3
$ws.Range("D3").Value = [double]"3.608291491697465E-10"
$ws.Range("E3").Value = [double]"3.608291491697465E-10"

# Row 4
$ws.Range("D4").Value = [double]"1.976246654820856E-29"
$ws.Range("E4").Value = [double]"1.976246654820856E-29"

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.01841013135184652
$ws.Range("E5").Value = 0.01841013135184652

# Row 6
$ws.Range("D6").Value = [double]"2.684844441776036E-44"
$ws.Range("E6").Value = [double]"2.684844441776036E-44"

# Row 7
$ws.Range("D7").Value = 0.9999999999999989
$ws.Range("E7").Value = [double]"1.110223024625157E-15"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.1189651221027286
$ws.Range("E8").Value = 0.8810348778972714

# Row 9
$ws.Range("D9").Value = 0.9999864915141926
$ws.Range("E9").Value = [double]"1.350848580738973E-05"

# Row 10
$ws.Range("D10").Value = [double]"9.843467562190021E-31"

# Row 11
$ws.Range("D11").Value = [double]"8.468257556382933E-88"
$ws.Range("F11").Value = 29.78525924682617
